$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content changes ---

# CRMUrl: point at the PROD CRM environment instead of UAT (keep original trailing
# run of non-breaking spaces + one regular space that padded the old value).
$nbsp = [char]0x00A0
$pad = ""
for ($i = 0; $i -lt 55; $i++) { $pad += $nbsp }
$ws.Range("B21").Value = "https://rpa-csc-prod.crm4.dynamics.com/main.aspx" + $pad + " "

# XcomUrl: point at the new intranet host instead of the old dev domino box
$ws.Range("B25").Value = "http://crosscompliancedatabase/intranet/xcompliance2015.nsf/"

# Mailbox values switched from raw mail addresses to the new SM- distribution
# list display names. Set in the order they first appear (SAG, RPA, then AH)
# so newly added shared strings line up the same way they do upstream.
$ws.Range("B42").Value = "SM-RPA-XC SAG Reports"
$ws.Range("B45").Value = "SM-RPA-XCRPAreports"
$ws.Range("B43").Value = "SM-RPA-XC Animal Health Stand alone and selected"

# Row 25 had an explicit 30pt row height; restore it to automatic/default height
$ws.Rows.Item(25).AutoFit()

# --- View state changes ---
$ws.Range("B45").Select()
